$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 5 de Octubre de 2020 a las 21:52"

# Row 4 - Estados Unidos
$ws.Range("B4").Value = 7664106
$ws.Range("C4").Value = 27194
$ws.Range("D4").Value = 4872671
$ws.Range("E4").Value = 2576583
$ws.Range("G4").Value = 241
$ws.Range("H4").Value = 214852

# Row 5 - India
$ws.Range("B5").Value = 6682073
$ws.Range("C5").Value = 59893
$ws.Range("D5").Value = 5659110
$ws.Range("E5").Value = 919363

# Row 26 - Alemania
$ws.Range("B26").Value = 304636
$ws.Range("C26").Value = 3065
$ws.Range("E26").Value = 31320
$ws.Range("G26").Value = 14
$ws.Range("H26").Value = 9616

# Row 111 - Haiti
$ws.Range("B111").Value = 8827
$ws.Range("C111").Value = 8
$ws.Range("E111").Value = 1606

# Row 120 - Malaui
$ws.Range("B120").Value = 5794
$ws.Range("C120").Value = 8
$ws.Range("E120").Value = 1073
$ws.Range("G120").Value = 1
$ws.Range("H120").Value = 180

# Row 131 - Ruanda
$ws.Range("B131").Value = 4867
$ws.Range("C131").Value = 1
$ws.Range("D131").Value = 3226
$ws.Range("E131").Value = 1612

# Row 134 - Bahamas
$ws.Range("B134").Value = 4452
$ws.Range("C134").Value = 43
$ws.Range("E134").Value = 1981
